$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new question code, new correct answer)
$changes = @{
    2  = @("M01A", 4)
    3  = @("M03A", 1)
    4  = @("M06A", 1)
    5  = @("M08A", 2)
    6  = @("O09A", 3)
    7  = @("M04B", 4)
    8  = @("M05B", 3)
    9  = @("M08B", 4)
    10 = @("O04B", 1)
    11 = @("O05B", 2)
    12 = @("M07C", 4)
    13 = @("M10C", 1)
    14 = @("O07C", 1)
    15 = @("O08C", 4)
    16 = @("O10C", 3)
    17 = @("N01D", 1)
    18 = @("N08D", 4)
    19 = @("N09D", 3)
    20 = @("O04D", 2)
    21 = @("O06D", 1)
    22 = @("M01E", 4)
    23 = @("M09E", 2)
    24 = @("O10E", 2)
    25 = @("N01F", 2)
    26 = @("N10F", 4)
    27 = @("O01F", 2)
    28 = @("M10G", 1)
    29 = @("O03G", 4)
    30 = @("O05G", 2)
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
}
